$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the grid of values (B1:D1, A2:D4) while keeping styles,
# and set A1 to the real published value (4).
$ws.Range("B1:D1").ClearContents()
$ws.Range("A2:D4").ClearContents()
$ws.Range("A1").Value = 4

# Reset the active selection back to A1.
$ws.Range("A1").Select() | Out-Null
